$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Fill in the "day 6" tracking column (L) for the 15-Nov sheet with the new
# readings captured on 20-Nov; dependent SUM/AVERAGE formulas recalc on save.
$ws1.Range("L8").Value = 4509
$ws1.Range("L9").Value = 1588
$ws1.Range("L10").Value = 7494
$ws1.Range("L11").Value = 29
$ws1.Range("L12").Value = 0
$ws1.Range("L13").Value = 324
$ws1.Range("L14").Value = 296
$ws1.Range("L15").Value = 47

# Move the selection/scroll position to where the latest update was made.
$ws1.Range("L15").Select()

# Duplicate the now-updated sheet (formulas + values) and drop the copy in
# right before the old, still-empty "Sheet2" - this copy becomes the new
# "Sheet1" that tracking continues on going forward.
$ws1.Copy($ws2)

# The original empty Sheet2 got pushed one slot to the right by the copy;
# remove it since the duplicate above now serves as the live Sheet1.
$oldSheet2 = $wb.Worksheets.Item(3)
$oldSheet2.Delete()

# Append a brand-new blank worksheet at the end to start tracking rice.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)

# Rename through unique temporary names first so the Sheet2/Sheet3 swap and
# the Sheet1 hand-off don't collide with names still in use mid-sequence.
$wb.Worksheets.Item(1).Name = "TMP1__"
$wb.Worksheets.Item(2).Name = "TMP2__"
$wb.Worksheets.Item(3).Name = "TMP3__"
$wb.Worksheets.Item(4).Name = "TMP4__"

$wb.Worksheets.Item(1).Name = "Sheet1_2(20Nov)"
$wb.Worksheets.Item(2).Name = "Sheet1"
$wb.Worksheets.Item(3).Name = "Sheet2"
$wb.Worksheets.Item(4).Name = "Sheet3"

$wb.Worksheets.Item(1).Activate()
